$wb = $excel.ActiveWorkbook

$updates = @(
  @{ Sheet="ALC"; Row=76; Sets=@{"H"=2890.9736; "I"=2893; "J"=2882; "K"=2893; "L"=2882; "M"=-2578; "N"=-3512}; Clears=@() },
  @{ Sheet="ALC"; Row=79; Sets=@{"H"=2890.9736; "I"=2893; "J"=2882; "K"=2893; "L"=2882; "M"=-1801; "N"=-5066}; Clears=@() },
  @{ Sheet="ALC"; Row=80; Sets=@{"H"=733.35297; "I"=649.3333; "J"=827.875; "K"=1947.9999; "L"=2483.625; "M"=-949.9999; "N"=-4479.625}; Clears=@() },
  @{ Sheet="ALC"; Row=83; Sets=@{"H"=733.35297; "I"=649.3333; "J"=827.875; "K"=5843.9997; "L"=7450.875; "M"=-851.9997000000003; "N"=-17434.875}; Clears=@() },
  @{ Sheet="ALC"; Row=88; Sets=@{"H"=1928; "I"=360.75; "J"=3495.25; "K"=360.75; "L"=3495.25; "M"=45.25; "N"=-4307.25}; Clears=@() },
  @{ Sheet="ALC"; Row=91; Sets=@{"H"=1928; "I"=360.75; "J"=3495.25; "K"=360.75; "L"=3495.25; "M"=1043.25; "N"=-6303.25}; Clears=@() },
  @{ Sheet="ALC"; Row=94; Sets=@{"H"=4578.6665; "I"=3225.8333; "J"=9990; "K"=3225.8333; "L"=9990; "M"=-2774.8333; "N"=-10892}; Clears=@() },
  @{ Sheet="ALC"; Row=96; Sets=@{"H"=526.4286; "I"=424.1; "J"=782.25; "K"=1272.3; "L"=2346.75; "M"=100.6999999999998; "N"=-5092.75}; Clears=@() },
  @{ Sheet="ALC"; Row=112; Sets=@{"H"=52632910; "I"=699.6; "J"=71430130; "K"=2098.8; "L"=214290390; "M"=-990.8000000000002; "N"=-214292606}; Clears=@() },
  @{ Sheet="ALC"; Row=127; Sets=@{"H"=808.3333; "I"=800; "J"=809.0909; "K"=2400; "L"=2427.2727; "M"=2560; "N"=-12347.2727}; Clears=@() },
  @{ Sheet="ALC"; Row=129; Sets=@{"H"=964.6087; "I"=323.75; "J"=1190.7941; "K"=971.25; "L"=3572.3823; "M"=4028.75; "N"=-13572.3823}; Clears=@() },
  @{ Sheet="ALC"; Row=132; Sets=@{"H"=3902.92; "I"=4220.5654; "J"=250; "K"=12661.6962; "L"=750; "M"=-10131.6962; "N"=-5810}; Clears=@() },
  @{ Sheet="ALC"; Row=137; Sets=@{"H"=26318670; "I"=1537.9259; "J"=90915260; "K"=4613.7777; "L"=272745780; "M"=-2063.7777; "N"=-272750880}; Clears=@() },
  @{ Sheet="ALC"; Row=138; Sets=@{"H"=2625.24; "I"=2052.543; "J"=2933.6155; "K"=6157.629000000001; "L"=8800.8465; "M"=-1017.629000000001; "N"=-19080.8465}; Clears=@() },
  @{ Sheet="ARM"; Row=32; Sets=@{"H"=21092.645; "I"=15392.138; "J"=103750; "K"=15392.138; "L"=103750; "M"=-15105.138; "N"=-104324}; Clears=@() },
  @{ Sheet="ARM"; Row=132; Sets=@{"H"=933930.7; "I"=1079520.5; "J"=2155.6; "K"=3238561.5; "L"=6466.799999999999; "M"=-3236031.5; "N"=-11526.8}; Clears=@() },
  @{ Sheet="BSM"; Row=20; Sets=@{"H"=3944; "I"=6166.6665; "J"=2277; "K"=6166.6665; "L"=2277; "M"=-5919.6665; "N"=-2771}; Clears=@() },
  @{ Sheet="BSM"; Row=25; Sets=@{"H"=3529.75; "I"=1373; "J"=10000; "K"=1373; "L"=10000; "M"=-1138}; Clears=@() },
  @{ Sheet="BSM"; Row=134; Sets=@{"H"=4652593.5; "I"=5129632; "J"=1472.5; "K"=15388896; "L"=4417.5; "M"=-15386361; "N"=-9487.5}; Clears=@() },
  @{ Sheet="CRP"; Row=7; Sets=@{"H"=52.833332; "I"=51.25; "J"=56; "K"=51.25; "L"=56; "M"=61.75; "N"=-282}; Clears=@() },
  @{ Sheet="CRP"; Row=22; Sets=@{"H"=212; "I"=200; "J"=350; "K"=200; "L"=350; "M"=150; "N"=-1050}; Clears=@() },
  @{ Sheet="CRP"; Row=31; Sets=@{"H"=20477.646; "I"=43452.168; "J"=7946.091; "K"=43452.168; "L"=7946.091; "M"=-43157.168; "N"=-8536.091}; Clears=@() },
  @{ Sheet="CRP"; Row=34; Sets=@{"H"=20477.646; "I"=43452.168; "J"=7946.091; "K"=43452.168; "L"=7946.091; "M"=-43250.168; "N"=-8350.091}; Clears=@() },
  @{ Sheet="CRP"; Row=88; Sets=@{"H"=0; "I"=0; "J"=0; "K"=0; "L"=0}; Clears=@("M","N") },
  @{ Sheet="CRP"; Row=91; Sets=@{"H"=0; "I"=0; "J"=0; "K"=0; "L"=0}; Clears=@("M","N") },
  @{ Sheet="CRP"; Row=99; Sets=@{"H"=101236.5; "I"=1281.375; "J"=501057; "K"=1281.375; "L"=501057; "M"=216.625; "N"=-504053}; Clears=@() },
  @{ Sheet="CRP"; Row=105; Sets=@{"H"=978; "I"=727.1429000000001; "J"=1563.3334; "K"=727.1429000000001; "L"=1563.3334; "M"=1019.8571; "N"=-5057.3334}; Clears=@() },
  @{ Sheet="CRP"; Row=126; Sets=@{"H"=101236.5; "I"=1281.375; "J"=501057; "K"=3844.125; "L"=1503171; "M"=-1374.125; "N"=-1508111}; Clears=@() },
  @{ Sheet="CUL"; Row=97; Sets=@{"H"=389.81818; "I"=525; "J"=312.57144; "K"=1575; "L"=937.71432; "M"=-1079; "N"=-1929.71432}; Clears=@() },
  @{ Sheet="CUL"; Row=98; Sets=@{"H"=346.66666; "I"=293.33334; "J"=400; "K"=880.0000200000001; "L"=1200; "M"=617.9999799999999; "N"=-4196}; Clears=@() },
  @{ Sheet="CUL"; Row=105; Sets=@{"H"=297335300; "I"=0; "J"=297335300; "K"=0; "L"=892005900; "N"=-892011142}; Clears=@("M") },
  @{ Sheet="CUL"; Row=122; Sets=@{"H"=6739832; "I"=15873680; "J"=745743.75; "K"=142863120; "L"=6711693.75; "M"=-142860670; "N"=-6716593.75}; Clears=@() },
  @{ Sheet="CUL"; Row=131; Sets=@{"H"=31535180; "I"=166681680; "J"=15153787; "K"=500045040; "L"=45461361; "M"=-500040000; "N"=-45471441}; Clears=@() },
  @{ Sheet="CUL"; Row=132; Sets=@{"H"=47620108; "I"=62500896; "J"=1595.8; "K"=562508064; "L"=14362.2; "M"=-562505534; "N"=-19422.2}; Clears=@() },
  @{ Sheet="GSM"; Row=107; Sets=@{"H"=337.2069; "I"=296.6111; "J"=403.63635; "K"=296.6111; "L"=403.63635; "M"=1623.3889; "N"=-4243.63635}; Clears=@() },
  @{ Sheet="GSM"; Row=113; Sets=@{"H"=1985.4286; "I"=2731.6667; "J"=1425.75; "K"=2731.6667; "L"=1425.75; "M"=-561.6667000000002; "N"=-5765.75}; Clears=@() },
  @{ Sheet="GSM"; Row=122; Sets=@{"H"=3927.5652; "I"=5011; "J"=2242.2222; "K"=15033; "L"=6726.6666; "M"=-12583; "N"=-11626.6666}; Clears=@() },
  @{ Sheet="GSM"; Row=132; Sets=@{"H"=1342.3158; "I"=985.6070999999999; "J"=2341.1; "K"=2956.8213; "L"=7023.299999999999; "M"=-426.8212999999996; "N"=-12083.3}; Clears=@() },
  @{ Sheet="WVR"; Row=126; Sets=@{"H"=2962.6875; "I"=2091.818; "J"=4878.6; "K"=6275.454000000001; "L"=14635.8; "M"=-3805.454000000001; "N"=-19575.8}; Clears=@() }
)

$colIndex = @{ "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "N"=14 }

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    foreach ($col in $u.Sets.Keys) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).Value = $u.Sets[$col]
    }
    foreach ($col in $u.Clears) {
        $ws.Cells.Item($u.Row, $colIndex[$col]).ClearContents()
    }
}